# Scheduled market-data refresh: update Leve profit/price figures per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 227.33333
$ws.Range("I2").Value = 233.71428
$ws.Range("J2").Value = 205
$ws.Range("K2").Value = 233.71428
$ws.Range("L2").Value = 205
$ws.Range("M2").Value = -120.71428
$ws.Range("N2").Value = -431
$ws.Range("H3").Value = 52500
$ws.Range("J3").Value = 52500
$ws.Range("L3").Value = 52500
$ws.Range("N3").Value = -52728
$ws.Range("H13").Value = 27499.25
$ws.Range("H76").Value = 6142.857
$ws.Range("I76").Value = 7250
$ws.Range("J76").Value = 4666.6665
$ws.Range("K76").Value = 7250
$ws.Range("L76").Value = 4666.6665
$ws.Range("M76").Value = -6935
$ws.Range("N76").Value = -5296.6665
$ws.Range("H79").Value = 6142.857
$ws.Range("I79").Value = 7250
$ws.Range("J79").Value = 4666.6665
$ws.Range("K79").Value = 7250
$ws.Range("L79").Value = 4666.6665
$ws.Range("M79").Value = -6158
$ws.Range("N79").Value = -6850.6665
$ws.Range("H102").Value = 52500
$ws.Range("J102").Value = 52500
$ws.Range("L102").Value = 52500
$ws.Range("N102").Value = -58990
$ws.Range("H107").Value = 334.15384
$ws.Range("I107").Value = 253.45454
$ws.Range("J107").Value = 778
$ws.Range("K107").Value = 253.45454
$ws.Range("L107").Value = 778
$ws.Range("M107").Value = 1666.54546
$ws.Range("N107").Value = -4618
$ws.Range("H112").Value = 2844.4
$ws.Range("J112").Value = 3305.5
$ws.Range("L112").Value = 9916.5
$ws.Range("N112").Value = -12132.5
$ws.Range("H116").Value = 13335567
$ws.Range("I116").Value = 40001900
$ws.Range("K116").Value = 40001900
$ws.Range("M116").Value = -39998458
$ws.Range("H129").Value = 1259.2
$ws.Range("J129").Value = 1428.3448
$ws.Range("L129").Value = 4285.0344
$ws.Range("N129").Value = -14285.0344
$ws.Range("H137").Value = 2818.7334
$ws.Range("I137").Value = 1632.1
$ws.Range("K137").Value = 4896.299999999999
$ws.Range("M137").Value = -2346.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1479.6842
$ws.Range("I74").Value = 1875
$ws.Range("J74").Value = 1040.4445
$ws.Range("K74").Value = 1875
$ws.Range("L74").Value = 1040.4445
$ws.Range("M74").Value = -1001
$ws.Range("N74").Value = -2788.4445
$ws.Range("H77").Value = 1479.6842
$ws.Range("I77").Value = 1875
$ws.Range("J77").Value = 1040.4445
$ws.Range("K77").Value = 9375
$ws.Range("L77").Value = 5202.2225
$ws.Range("M77").Value = -5007
$ws.Range("N77").Value = -13938.2225
$ws.Range("H88").Value = 3022.2222
$ws.Range("I88").Value = 2600
$ws.Range("K88").Value = 2600
$ws.Range("M88").Value = -2194
$ws.Range("H91").Value = 3022.2222
$ws.Range("I91").Value = 2600
$ws.Range("K91").Value = 2600
$ws.Range("M91").Value = -1196

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3161.12
$ws.Range("I105").Value = 3193.4
$ws.Range("J105").Value = 3112.7
$ws.Range("K105").Value = 3193.4
$ws.Range("L105").Value = 3112.7
$ws.Range("M105").Value = -1446.4
$ws.Range("N105").Value = -6606.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 442.75
$ws.Range("I22").Value = 241.71428
$ws.Range("J22").Value = 1850
$ws.Range("K22").Value = 241.71428
$ws.Range("L22").Value = 1850
$ws.Range("M22").Value = 108.28572
$ws.Range("N22").Value = -2550
$ws.Range("H31").Value = 2295.2122
$ws.Range("I31").Value = 1418.3572
$ws.Range("J31").Value = 7205.6
$ws.Range("K31").Value = 1418.3572
$ws.Range("L31").Value = 7205.6
$ws.Range("M31").Value = -1123.3572
$ws.Range("N31").Value = -7795.6
$ws.Range("H34").Value = 2295.2122
$ws.Range("I34").Value = 1418.3572
$ws.Range("J34").Value = 7205.6
$ws.Range("K34").Value = 1418.3572
$ws.Range("L34").Value = 7205.6
$ws.Range("M34").Value = -1216.3572
$ws.Range("N34").Value = -7609.6
$ws.Range("H94").Value = 751.4211
$ws.Range("I94").Value = 866.3333
$ws.Range("J94").Value = 698.38464
$ws.Range("K94").Value = 866.3333
$ws.Range("L94").Value = 698.38464
$ws.Range("M94").Value = -415.3333
$ws.Range("N94").Value = -1600.38464

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2244.5789
$ws.Range("J5").Value = 1980
$ws.Range("L5").Value = 5940
$ws.Range("N5").Value = -6164
$ws.Range("H121").Value = 36995.93
$ws.Range("J121").Value = 63890.375
$ws.Range("L121").Value = 191671.125
$ws.Range("N121").Value = -194291.125
$ws.Range("H131").Value = 12659418
$ws.Range("J131").Value = 12988189
$ws.Range("L131").Value = 38964567
$ws.Range("N131").Value = -38974647
$ws.Range("H135").Value = 2244.5789
$ws.Range("J135").Value = 1980
$ws.Range("L135").Value = 17820
$ws.Range("N135").Value = -22890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 30718.75
$ws.Range("J93").Value = 30718.75
$ws.Range("L93").Value = 30718.75
$ws.Range("N93").Value = -34462.75
$ws.Range("H102").Value = 2720.2
$ws.Range("I102").Value = 2291.4443
$ws.Range("K102").Value = 2291.4443
$ws.Range("M102").Value = -669.4443000000001
$ws.Range("H109").Value = 16951
$ws.Range("J109").Value = 16951
$ws.Range("L109").Value = 16951
$ws.Range("N109").Value = -19031
$ws.Range("H113").Value = 3170.0557
$ws.Range("I113").Value = 1635.1666
$ws.Range("J113").Value = 3937.5
$ws.Range("K113").Value = 1635.1666
$ws.Range("L113").Value = 3937.5
$ws.Range("M113").Value = 534.8334
$ws.Range("N113").Value = -8277.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6346.5
$ws.Range("J40").Value = 4699.75
$ws.Range("L40").Value = 4699.75
$ws.Range("N40").Value = -4971.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 11666
$ws.Range("J15").Value = 11666
$ws.Range("L15").Value = 11666
$ws.Range("N15").Value = -12242
$ws.Range("H64").Value = 31371.334
$ws.Range("J64").Value = 31371.334
$ws.Range("L64").Value = 31371.334
$ws.Range("N64").Value = -31867.334
$ws.Range("H67").Value = 31371.334
$ws.Range("J67").Value = 31371.334
$ws.Range("L67").Value = 31371.334
$ws.Range("N67").Value = -33087.334
$ws.Range("H81").Value = 123044.555
$ws.Range("I81").Value = 123044.555
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 246089.11
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -245028.11
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 123044.555
$ws.Range("I84").Value = 123044.555
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 1230445.55
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -1225141.55
$ws.Range("N84").ClearContents()
$ws.Range("H95").Value = 139750
$ws.Range("J95").Value = 139750
$ws.Range("L95").Value = 139750
$ws.Range("N95").Value = -145242
$ws.Range("H107").Value = 816.8
$ws.Range("I107").Value = 794.125
$ws.Range("J107").Value = 907.5
$ws.Range("K107").Value = 2382.375
$ws.Range("L107").Value = 2722.5
$ws.Range("M107").Value = -462.375
$ws.Range("N107").Value = -6562.5
$ws.Range("H113").Value = 366.65
$ws.Range("I113").Value = 429.625
$ws.Range("J113").Value = 324.66666
$ws.Range("K113").Value = 1288.875
$ws.Range("L113").Value = 973.9999799999999
$ws.Range("M113").Value = 881.125
$ws.Range("N113").Value = -5313.99998
